$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text would otherwise be auto-parsed by Excel as a number
# (e.g. "1.00" -> 1, losing the trailing zero / textual form).
# Temporarily mark them as Text, assign the literal string, then restore the
# default "Normal" cell style so no stray formatting is left behind.
function Set-TextValue($range, $text) {
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.Style = "Normal"
}

$ws.Range("D2").Value = "60.866.41"
$ws.Range("E2").Value = "  -3.29%  "

$ws.Range("D3").Value = "2.908.17"
$ws.Range("E3").Value = "  -3.92%  "

Set-TextValue $ws.Range("D4") "1.00"
$ws.Range("E4").Value = "  -0.03%  "

Set-TextValue $ws.Range("D5") "590.24"
$ws.Range("E5").Value = "  -0.91%  "

Set-TextValue $ws.Range("D6") "144.75"
$ws.Range("E6").Value = "  -5.37%  "

$ws.Range("E7").Value = "  +0.05%  "

$ws.Range("E8").Value = "  -1.53%  "

$ws.Range("D9").Value = "2.907.80"
$ws.Range("E9").Value = "  -3.83%  "

$ws.Range("E10").Value = "  -4.18%  "

Set-TextValue $ws.Range("D11") "0.144"
$ws.Range("E11").Value = "  -4.00%  "

$ws.Range("E12").Value = "  -4.11%  "

$ws.Range("E13").Value = "  -2.68%  "

Set-TextValue $ws.Range("D14") "33.50"
$ws.Range("E14").Value = "  -6.07%  "

Set-TextValue $ws.Range("D15") "0.126"
$ws.Range("E15").Value = "  +0.38%  "

$ws.Range("D16").Value = "3.389.75"
$ws.Range("E16").Value = "  -3.92%  "

$ws.Range("D17").Value = "60.796.63"
$ws.Range("E17").Value = "  -3.38%  "

Set-TextValue $ws.Range("D18") "6.72"
$ws.Range("E18").Value = "  -4.93%  "

$ws.Range("D19").Value = "2.904.39"
$ws.Range("E19").Value = "  -3.97%  "

Set-TextValue $ws.Range("D20") "430.72"
$ws.Range("E20").Value = "  -4.09%  "

$ws.Range("E21").Value = "  -4.92%  "

$ws.Range("E22").Value = "  -1.87%  "

$ws.Range("E23").Value = "  -5.74%  "

Set-TextValue $ws.Range("D24") "81.63"
$ws.Range("E24").Value = "  -1.48%  "

Set-TextValue $ws.Range("D25") "10.86"
$ws.Range("E25").Value = "  -4.90%  "

$ws.Range("E26").Value = "  -3.57%  "

Set-TextValue $ws.Range("D27") "12.02"
$ws.Range("E27").Value = "  -2.61%  "

$ws.Range("E28").Value = "  +0.04%  "

$ws.Range("E29").Value = "  +0.50%  "

Set-TextValue $ws.Range("D30") "1.00"
$ws.Range("E30").Value = "  -0.06%  "

$ws.Range("E31").Value = "  -2.61%  "

$ws.Range("E32").Value = "  -5.57%  "

Set-TextValue $ws.Range("D33") "26.60"

$ws.Range("E34").Value = "  -2.87%  "

$ws.Range("E35").Value = "  -2.07%  "

$ws.Range("E36").Value = "  -3.33%  "

$ws.Range("E37").Value = "  -4.67%  "

Set-TextValue $ws.Range("D38") "3.00"
$ws.Range("E38").Value = "  -4.00%  "

Set-TextValue $ws.Range("D39") "49.62"
$ws.Range("E39").Value = "  -1.83%  "

Set-TextValue $ws.Range("D40") "0.124"
$ws.Range("E40").Value = "  -4.53%  "

$ws.Range("E41").Value = "  -4.57%  "

Set-TextValue $ws.Range("D42") "8.62"
$ws.Range("E42").Value = "  -4.51%  "

Set-TextValue $ws.Range("D43") "0.291"
$ws.Range("E43").Value = "  -4.88%  "

Set-TextValue $ws.Range("D44") "40.06"
$ws.Range("E44").Value = "  -9.97%  "

$ws.Range("E45").Value = "  -3.05%  "

Set-TextValue $ws.Range("D46") "374.31"
$ws.Range("E46").Value = "  -4.30%  "

$ws.Range("D47").Value = "2.703.08"
$ws.Range("E47").Value = "  -0.19%  "

Set-TextValue $ws.Range("D48") "130.29"
$ws.Range("E48").Value = "  -2.70%  "

$ws.Range("E50").Value = "  -10.14%  "

$ws.Range("E51").Value = "  -2.18%  "
